$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells that hold plain numeric-looking strings
# to be written as TEXT (matching the original inlineStr/shared-string cells)
# instead of being auto-converted to numbers by Excel, then restore the
# default "Normal" cell style so no stray formatting is introduced.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '62.191.55'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.73%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.185.23'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -3.68%  '

$ws.Range("E4").Value = '  -0.01%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '593.32'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.16%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '135.69'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -5.95%  '

$ws.Range("E7").Value = '  -0.02%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '3.182.12'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.71%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.506'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.21%  '

$ws.Range("E10").Value = '  -4.22%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '5.31'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -2.99%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.452'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -4.15%  '

$ws.Range("E13").Value = '  -5.09%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '33.49'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -4.13%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.714.92'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -3.53%  '

$ws.Range("E16").Value = '  -0.29%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.191.14'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -3.45%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '62.349.25'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.63%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '6.66'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.40%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '462.42'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.99%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '13.98'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.29%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.709'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -4.42%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.67'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -4.29%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '13.41'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.18%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '83.69'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("E26").Value = '  -0.18%  '

$ws.Range("E27").Value = '  -2.63%  '

$ws.Range("E28").Value = '  -0.04%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '7.92'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -4.18%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '6.90'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -5.22%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.06'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -4.18%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '27.28'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.03%  '

$ws.Range("E33").Value = '  -4.38%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '2.42'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -5.34%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.04'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -5.78%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '5.85'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.53%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '51.50'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -3.46%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0691'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -8.98%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.0389'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.003.58'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.99%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '412.33'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -4.42%  '

$ws.Range("E42").Value = '  +4.27%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '8.08'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -4.35%  '

$ws.Range("E44").Value = '  -6.37%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.252'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -6.77%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.16'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -3.29%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '36.04'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.62%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '25.80'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.61%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '124.08'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.04%  '
